# Correct capitalization in project title
#
# "Evaluation of tool support to validate the Object Calisthenics"
#   -> "Evaluation and Prototypical Implementation of Tool Support to
#       Validate the Object Calisthenics"
#
# Word drops its internal "last edit" bookmark (_GoBack) at the spot of the
# most recent change, so after re-typing the title we relocate that
# bookmark from its old position (in the "Usually the developer..."
# paragraph) onto the title paragraph.

$d = $word.ActiveDocument

$titleParagraph = $d.Paragraphs(1)

# Rewrite the title text itself.
$d.Content.Find.Execute(
    "Evaluation of tool support to validate the Object Calisthenics",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Evaluation and Prototypical Implementation of Tool Support to Validate the Object Calisthenics",
    2)

# The old "_GoBack" bookmark sat mid-word, splitting "softwar" / "e with
# little requirements..." into two runs. Re-typing that span as one
# contiguous run both heals the split and removes the stale bookmark that
# was anchored inside it.
$d.Content.Find.Execute(
    "minimalistic software with little requirements",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "minimalistic software with little requirements",
    2)

# Move the "_GoBack" bookmark (re-adding with the same name relocates it,
# since a document can only have one bookmark per name) so it wraps the
# paragraph that was just edited, same as Word does automatically.
$d.Bookmarks.Add("_GoBack", $titleParagraph.Range)
